$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("you", "あなた"),
    @("chair", "いす"),
    @("dog", "犬|いぬ"),
    @("souvenir", "お土産|おみやげ"),
    @("child", "子供|こども"),
    @("rice; meal", "ご飯|ごはん"),
    @("picture; photograph", "写真|しゃしん"),
    @("desk", "机|つくえ"),
    @("letter", "手紙|てがみ"),
    @("cat", "猫|ねこ"),
    @("bread", "パン"),
    @("person", "人|ひと"),
    @("e-mail", "メール"),
    @("the day before yesterday", "おととい"),
    @("the day after tomorrow", "あさって"),
    @("two weeks ago", "二週間前|にしゅうかんまえ"),
    @("this week", "今週|こんしゅう"),
    @("next week", "来週|らいしゅう"),
    @("the week after next", "再来週|さらいしゅう"),
    @("two months ago", "二か月前|にかげつまえ"),
    @("last month", "先月|せんげつ"),
    @("this month", "今月|こんげつ"),
    @("next month", "来月|らいげつ"),
    @("the month after next", "再来月|さらいげつ"),
    @("the year before last", "おととし"),
    @("last year", "去年|きょねん"),
    @("this year", "今年|ことし"),
    @("next year", "来年|らいねん"),
    @("the year after next", "再来年|さらいねん"),
    @("part-time job", "アルバイト"),
    @("shopping", "買い物|かいもの"),
    @("class", "クラス"),
    @("temple", "お寺|おてら"),
    @("park", "公園|こうえん"),
    @("supermarket", "スーパー"),
    @("department store", "デパート"),
    @("bus stop", "バス停|バスてい"),
    @("hospital", "病院|びょういん"),
    @("hotel", "ホテル"),
    @("bookstore", "本屋|ほんや"),
    @("town; city", "町|まち"),
    @("restaurant", "レストラン"),
    @("yesterday", "昨日|きのう"),
    @("...hours", "～時間|～じかん"),
    @("one hour", "一時間|いちじかん"),
    @("last week", "先週|せんしゅう"),
    @("when...; at the time of...", "時|とき"),
    @("Monday", "月曜日|げつようび"),
    @("Tuesday", "火曜日|かようび"),
    @("Wednesday", "水曜日|すいようび"),
    @("Thursday", "木曜日|もくようび"),
    @("Friday", "金曜日|きんようび"),
    @("to meet; to see (a person)", "会う|あう"),
    @("there is...", "ある"),
    @("to buy", "買う|かう"),
    @("to write", "書く|かく"),
    @("to take (a picture)", "撮る|とる"),
    @("to wait", "待つ|まつ"),
    @("to understand", "わかる"),
    @("(a person) is in...; stays at...", "いる"),
    @("about (approximate measurement)", "～ぐらい"),
    @("I'm sorry.", "ごめんなさい"),
    @("so; therefore", "だから"),
    @("many; a lot", "たくさん"),
    @("together with (a person)", "～と"),
    @("why", "どうして"),
    @("alone", "一人で|ひとりで"),
    @("right", "右|みぎ"),
    @("left", "左|ひだり"),
    @("front", "前|まえ"),
    @("back", "後ろ|うしろ"),
    @("inside", "中|なか"),
    @("on", "上|うえ"),
    @("under", "下|した"),
    @("near; nearby", "近く|ちかく"),
    @("next", "隣|となり"),
    @("between", "間|あいだ"),
    @("1", "ついたち"),
    @("2", "ふつか"),
    @("3", "みっか"),
    @("4", "よっか"),
    @("5", "いつか"),
    @("6", "むいか"),
    @("7", "なのか"),
    @("8", "ようか"),
    @("9", "ここのか"),
    @("10", "とおか"),
    @("11", "じゅういちにち"),
    @("12", "じゅうににち"),
    @("13", "じゅうさんにち"),
    @("14", "じゅうよっか"),
    @("15", "じゅうごにち"),
    @("16", "じゅうろくにち"),
    @("17", "じゅうしちにち"),
    @("18", "じゅうはちにち"),
    @("19", "じゅうくにち"),
    @("20", "はつか"),
    @("21", "にじゅういちにち"),
    @("22", "にじゅうににち"),
    @("23", "にじゅうさんにち"),
    @("24", "にじゅうよっか"),
    @("25", "にじゅうごにち"),
    @("26", "にじゅうろくにち"),
    @("27", "にじゅうしちにち"),
    @("28", "にじゅうはちにち"),
    @("29", "にじゅうくにち"),
    @("30", "さんじゅうにち"),
    @("31", "さんじゅういちにち"),
    @("January", "一月|いちがつ"),
    @("February", "二月|にがつ"),
    @("March", "三月|さんがつ"),
    @("April", "四月|しがつ"),
    @("May", "五月|ごがつ"),
    @("June", "六月|ろくがつ"),
    @("July", "七月|しちがつ"),
    @("August", "八月|はちがつ"),
    @("September", "九月|くがつ"),
    @("October", "十月|じゅうがつ"),
    @("November", "十一月|じゅういちがつ"),
    @("December", "十二月|じゅうにがつ"),
    @("Japan", "日本|にほん"),
    @("Sunday", "日曜日|にちようび"),
    @("every day", "毎日|まいにち"),
    @("Mother's day", "母の日|ははのひ"),
    @("diary", "日記|にっき"),
    @("three days", "三日|みっか"),
    @("book", "本|ほん"),
    @("Japanese language", "日本語|にほんご"),
    @("Mr./Ms. Yamamoto", "山本さん|やまもとさん"),
    @("Japanese people", "日本人|にほんじん"),
    @("alone", "一人で|ひとりで"),
    @("this person", "この人|このひと"),
    @("three people", "三人|さんにん"),
    @("Monday", "月曜日|げつようび"),
    @("January", "一月|いちがつ"),
    @("moon", "月|つき"),
    @("this month", "今月|こんげつ"),
    @("one month", "一か月|いっかげつ"),
    @("Tuesday", "火曜日|かようび"),
    @("fire", "火|ひ"),
    @("volcano", "火山|かざん"),
    @("fireworks", "花火|はなび"),
    @("Mars", "火星|かせい"),
    @("Wednesday", "水曜日|すいようび"),
    @("water", "水|みず"),
    @("swimming", "水泳|すいえい"),
    @("water supply", "水道|すいどう"),
    @("bathing suit", "水着|みずぎ"),
    @("Thursday", "木曜日|もくようび"),
    @("tree", "木|き"),
    @("Mr./Ms. Kimura", "木村さん|きむらさん"),
    @("Friday", "金曜日|きんようび"),
    @("money", "お金|おかね"),
    @("charge", "料金|りょうきん"),
    @("rich person", "お金持ち|おかねもち"),
    @("Saturday", "土曜日|どようび"),
    @("soil", "土|つち"),
    @("land", "土地|とち"),
    @("clay", "粘土|ねんど"),
    @("Sunday", "日曜日|にちようび"),
    @("day of the week", "曜日|ようび"),
    @("top; above", "上|うえ"),
    @("good at", "上手な|じょうずな"),
    @("rooftop", "屋上|おくじょう"),
    @("to go up", "上る|のぼる"),
    @("under", "下|した"),
    @("subway", "地下鉄|ちかてつ"),
    @("poor at", "下手な|へたな"),
    @("Please give/do...", "下さい|ください"),
    @("inside", "中|なか"),
    @("China", "中国|ちゅうごく"),
    @("junior high school", "中学|ちゅうがく"),
    @("all year round", "一年中|いちねんじゅう"),
    @("half past three", "三時半|さんじはん"),
    @("half", "半分|はんぶん"),
    @("half a year", "半年|はんとし"),
    @("half price", "半額|はんがく"),
    @("Sunday", "日曜日"),
    @("Monday", "月曜日"),
    @("Tuesday", "火曜日"),
    @("Wednesday", "水曜日"),
    @("Thursday", "木曜日"),
    @("Friday", "金曜日"),
    @("Saturday", "土曜日")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $data[$i][0]
    $ws.Cells.Item($i + 2, 2).Value = $data[$i][1]
}
